# "Generate Report for Archive" - refresh the localization-status report:
#   * cells that still show "Ready for handoff" move on to "In Translation"
#   * the status columns are re-sized (narrower) to fit the new text

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status is mirrored per-locale in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRange = $wsOverview.Range("E2:F3")
foreach ($cell in $overviewRange.Cells) {
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: status lives in column C ---
$localeSheets = "zh-cn", "de-de"
foreach ($sheetName in $localeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $statusRange = $ws.Range("C2:C3")
    foreach ($cell in $statusRange.Cells) {
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
